# Add data for 2022-12-08
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet and update its "through" date references
$ws.Name = "Through 2022-11-30"
$ws.Range("A12").Value = "November (through 11-30)"

# Update November row (row 12) with the new day's counts
$ws.Range("B12").Value = 33
$ws.Range("C12").Value = 77
$ws.Range("D12").Value = 111
$ws.Range("F12").Value = 52
$ws.Range("G12").Value = 207
$ws.Range("H12").Value = 202
$ws.Range("I12").Value = 118

# Update Total row (row 13) to reflect the new counts
$ws.Range("B13").Value = 291
$ws.Range("C13").Value = 563
$ws.Range("D13").Value = 821
$ws.Range("F13").Value = 534
$ws.Range("G13").Value = 1264
$ws.Range("H13").Value = 1643
$ws.Range("I13").Value = 1516
